$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AY1) onto the
# three new header cells so they pick up the bold/centered/bordered style.
$ws.Range("AY1").Copy()
$ws.Range("AZ1:BB1").PasteSpecial(-4122)

# New header labels (appended as new shared-string entries / columns)
$ws.Range("AZ1").Value = "m_adult_literacy_pct"
$ws.Range("BA1").Value = "m_homicides_per_100k"
$ws.Range("BB1").Value = "m_tax_revenue_pct_gdp"

# New column data for each data row (rows 2-15)
$ws.Range("AZ2").Value = 0
$ws.Range("BA2").Value = 0
$ws.Range("BB2").Value = 1

$ws.Range("AZ3").Value = 1
$ws.Range("BA3").Value = 1
$ws.Range("BB3").Value = 1

$ws.Range("AZ4").Value = 1
$ws.Range("BA4").Value = 1
$ws.Range("BB4").Value = 1

$ws.Range("AZ5").Value = 1
$ws.Range("BA5").Value = 1
$ws.Range("BB5").Value = 0

$ws.Range("AZ6").Value = 1
$ws.Range("BA6").Value = 1
$ws.Range("BB6").Value = 0

$ws.Range("AZ7").Value = 1
$ws.Range("BA7").Value = 0
$ws.Range("BB7").Value = 1

$ws.Range("AZ8").Value = 1
$ws.Range("BA8").Value = 1
$ws.Range("BB8").Value = 0

$ws.Range("AZ9").Value = 1
$ws.Range("BA9").Value = 0
$ws.Range("BB9").Value = 1

$ws.Range("AZ10").Value = 1
$ws.Range("BA10").Value = 1
$ws.Range("BB10").Value = 1

$ws.Range("AZ11").Value = 0
$ws.Range("BA11").Value = 1
$ws.Range("BB11").Value = 0

$ws.Range("AZ12").Value = 1
$ws.Range("BA12").Value = 1
$ws.Range("BB12").Value = 1

$ws.Range("AZ13").Value = 0
$ws.Range("BA13").Value = 0
$ws.Range("BB13").Value = 0

$ws.Range("AZ14").Value = 1
$ws.Range("BA14").Value = 1
$ws.Range("BB14").Value = 1

$ws.Range("AZ15").Value = 1
$ws.Range("BA15").Value = 1
$ws.Range("BB15").Value = 1

Write-Output "edit complete"
